$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''70.312.63'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '''3.618.52'
$ws.Range('E3').Value = '  +2.43%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''602.41'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').Value = '''195.58'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.213'
$ws.Range('E9').Value = '  +4.57%  '
$ws.Range('D10').Value = '''0.646'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').Value = '''53.25'
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '''9.57'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').Value = '''4.188.81'
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('D15').Value = '''601.00'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '''12.97'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = '''70.437.18'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '''3.617.41'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('D19').Value = '''19.07'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = '''18.62'
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('D23').Value = '''5.21'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').Value = '''103.01'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').Value = '''9.72'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').Value = '''33.81'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').Value = '''4.73'
$ws.Range('E30').Value = '  +9.20%  '
$ws.Range('D31').Value = '''7.31'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').Value = '''12.29'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('D34').Value = '''63.32'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = '''0.0₃0884'
$ws.Range('E35').Value = '  +2.77%  '
$ws.Range('D36').Value = '''3.936.28'
$ws.Range('E36').Value = '  +5.19%  '
$ws.Range('D37').Value = '''531.73'
$ws.Range('E37').Value = '  +9.50%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '''3.05'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  -2.19%  '
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('D45').Value = '''3.61'
$ws.Range('E45').Value = '  +9.04%  '
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('D50').Value = '''0.000250'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('E51').Value = '  +1.66%  '
